$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 21.84976866666667
$ws.Range("H2").Value = 65.549306
$ws.Range("I2").Value = 0.05020018890879543
$ws.Range("J2").Value = 0.05020018890879543
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 1.398034
$ws.Range("N2").Value = 4.194102
$ws.Range("O2").Value = 0.139066772576779
$ws.Range("P2").Value = 0.139066772576779
$ws.Range("Q2").Value = 30.54671948813467
$ws.Range("R2").Value = 274.920475393212
$ws.Range("S2").Value = 0.006981178254290795
$ws.Range("T2").Value = 0.006981178254290797

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 21.84976866666667
$ws.Range("H3").Value = 65.549306
$ws.Range("I3").Value = 0.05020018890879543
$ws.Range("J3").Value = 0.05020018890879543
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 7.939250333333333
$ws.Range("N3").Value = 23.817751
$ws.Range("O3").Value = 0.7897418235434783
$ws.Range("P3").Value = 0.7897418235434784
$ws.Range("Q3").Value = 173.4707831700895
$ws.Range("R3").Value = 1561.237048530806
$ws.Range("S3").Value = 0.0396451887310592
$ws.Range("T3").Value = 0.0396451887310592

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 21.84976866666667
$ws.Range("H4").Value = 65.549306
$ws.Range("I4").Value = 0.05020018890879543
$ws.Range("J4").Value = 0.05020018890879543
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.715685
$ws.Range("N4").Value = 2.147055
$ws.Range("O4").Value = 0.0711914038797426
$ws.Range("P4").Value = 0.0711914038797426
$ws.Range("Q4").Value = 15.63755168820333
$ws.Range("R4").Value = 140.73796519383
$ws.Range("S4").Value = 0.00357382192344543
$ws.Range("T4").Value = 0.00357382192344543

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 385.0524703333334
$ws.Range("H5").Value = 1155.157411
$ws.Range("I5").Value = 0.8846641374295412
$ws.Range("J5").Value = 0.8846641374295412
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 1.398034
$ws.Range("N5").Value = 4.194102
$ws.Range("O5").Value = 0.139066772576779
$ws.Range("P5").Value = 0.139066772576779
$ws.Range("Q5").Value = 538.3164453099914
$ws.Range("R5").Value = 4844.848007789923
$ws.Range("S5").Value = 0.1230273864067463
$ws.Range("T5").Value = 0.1230273864067464

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 385.0524703333334
$ws.Range("H6").Value = 1155.157411
$ws.Range("I6").Value = 0.8846641374295412
$ws.Range("J6").Value = 0.8846641374295412
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 7.939250333333333
$ws.Range("N6").Value = 23.817751
$ws.Range("O6").Value = 0.7897418235434783
$ws.Range("P6").Value = 0.7897418235434784
$ws.Range("Q6").Value = 3057.02795344474
$ws.Range("R6").Value = 27513.25158100266
$ws.Range("S6").Value = 0.6986562691171242
$ws.Range("T6").Value = 0.6986562691171243

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 385.0524703333334
$ws.Range("H7").Value = 1155.157411
$ws.Range("I7").Value = 0.8846641374295412
$ws.Range("J7").Value = 0.8846641374295412
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 0.715685
$ws.Range("N7").Value = 2.147055
$ws.Range("O7").Value = 0.0711914038797426
$ws.Range("P7").Value = 0.0711914038797426
$ws.Range("Q7").Value = 275.5762772305117
$ws.Range("R7").Value = 2480.186495074605
$ws.Range("S7").Value = 0.06298048190567059
$ws.Range("T7").Value = 0.06298048190567059

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 28.350479
$ws.Range("H8").Value = 85.05143699999999
$ws.Range("I8").Value = 0.06513567366166337
$ws.Range("J8").Value = 0.06513567366166337
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 1.398034
$ws.Range("N8").Value = 4.194102
$ws.Range("O8").Value = 0.139066772576779
$ws.Range("P8").Value = 0.139066772576779
$ws.Range("Q8").Value = 39.63493355828599
$ws.Range("R8").Value = 356.714402024574
$ws.Range("S8").Value = 0.009058207915741832
$ws.Range("T8").Value = 0.009058207915741833

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 28.350479
$ws.Range("H9").Value = 85.05143699999999
$ws.Range("I9").Value = 0.06513567366166337
$ws.Range("J9").Value = 0.06513567366166337
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 7.939250333333333
$ws.Range("N9").Value = 23.817751
$ws.Range("O9").Value = 0.7897418235434783
$ws.Range("P9").Value = 0.7897418235434784
$ws.Range("Q9").Value = 225.0815498509096
$ws.Range("R9").Value = 2025.733948658187
$ws.Range("S9").Value = 0.05144036569529494
$ws.Range("T9").Value = 0.05144036569529495

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 28.350479
$ws.Range("H10").Value = 85.05143699999999
$ws.Range("I10").Value = 0.06513567366166337
$ws.Range("J10").Value = 0.06513567366166337
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 0.715685
$ws.Range("N10").Value = 2.147055
$ws.Range("O10").Value = 0.0711914038797426
$ws.Range("P10").Value = 0.0711914038797426
$ws.Range("Q10").Value = 20.290012563115
$ws.Range("R10").Value = 182.610113068035
$ws.Range("S10").Value = 0.00463710005062659
$ws.Range("T10").Value = 0.00463710005062659

